$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A14").Value = 37
$ws.Range("B14").Value = "latest timimg updated"
$ws.Range("C14").Value = "riya-morankar"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "edit1 to main"
$ws.Range("F14").NumberFormat = "@"
$ws.Range("F14").Value = "2025-06-18"
